$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows that get marked as completed on Milestone III (E = "III", F = "X")
$rows = @(29, 33, 40, 52, 64, 66, 67)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "III"
    $ws.Range("F$r").Value = "X"
}

# New source citations added to column A, rows 97-103
$ws.Range("A97").Value = "Unity Quaternion Math: https://docs.unity3d.com/ScriptReference/Quaternion.html"
$ws.Range("A98").Value = "Tutorial - Compute Shader Filters http://www.codinglabs.net/tutorial_compute_shaders_filters.aspx"
$ws.Range("A99").Value = "Cascaded Shadow Maps https://docs.microsoft.com/en-us/windows/desktop/dxtecharts/cascaded-shadow-maps"
$ws.Range("A100").Value = "OpenGL Physically Based Rendering : https://learnopengl.com/PBR/Theory"
$ws.Range("A101").Value = "OpenGL Bloom :  https://learnopengl.com/Advanced-Lighting/Bloom"
$ws.Range("A102").Value = "Sky HDRI Textures : https://hdrihaven.com/"
$ws.Range("A103").Value = "Shadow Filtering : https://docs.cryengine.com/display/SDKDOC4/Shadows+in+CryENGINE"

# Update sheet view scroll/selection state
$ws.Application.ActiveWindow.ScrollRow = 85
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A105").Select()

$wb.Save()
